$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table of worker account-statement rows (Tipo Doc, N Doc, Nombre, Periodo, Valor Mora, Salario Basico)
# for the range B16:G89. Previous (stale) statement data is fully replaced: Edgardo Fernandez
# Morales now has a single updated period (1608) in row 16, and Fermin Segundo Hernandez Ocando's
# 73 overdue periods (1607-2207) are listed in descending order in rows 17-89.
$rows = @(
    @("CC","73150584","EDGARDO FERNANDEZ MORALES","1608",57836,2040690),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2207",183333,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2206",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2205",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2204",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2203",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2202",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2201",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2112",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2111",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2110",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2109",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2108",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2107",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2106",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2105",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2104",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2103",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2102",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2101",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2012",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2011",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2010",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2009",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2008",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2007",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2006",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2005",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2004",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2003",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2002",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","2001",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1912",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1911",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1910",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1909",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1908",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1907",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1906",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1905",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1904",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1903",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1902",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1901",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1812",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1811",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1810",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1809",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1808",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1807",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1806",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1805",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1804",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1803",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1802",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1801",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1712",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1711",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1710",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1709",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1708",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1707",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1706",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1705",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1704",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1703",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1702",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1701",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1612",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1611",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1610",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1609",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1608",220000,5500000),
    @("CE","394833","FERMIN SEGUNDO HERNANDEZ OCANDO","1607",220000,5500000)
)

$n = $rows.Count
$cols = 6
$arr = New-Object "object[,]" $n, $cols
for ($i = 0; $i -lt $n; $i++) {
    for ($j = 0; $j -lt $cols; $j++) {
        $arr[$i, $j] = $rows[$i][$j]
    }
}

$ws.Range("B16:G89").Value = $arr
